$d = $word.ActiveDocument

# The document's headers/footers each carry one inline picture (a Pearson
# logo PNG in the footers, a BTEC logo JPG in the headers). Every copy of
# the PNG was named "image2.png" and every copy of the JPG was named
# "image1.jpg". The edit swaps those display/part names around:
#   footers: image2.png -> image1.png
#   headers: image1.jpg -> image2.jpg
# InlineShape has no writable "Name" in the real Word object model, so we
# briefly convert to a floating Shape (which does expose .Name), rename it,
# then convert back to an inline picture so the layout is unchanged.

function Rename-InlineShapes($story, $oldName, $newName) {
    for ($k = 1; $k -le $story.Range.InlineShapes.Count; $k++) {
        $ishp = $story.Range.InlineShapes.Item($k)
        $shp = $ishp.ConvertToShape()
        if ($shp.Name -eq $oldName) {
            $shp.Name = $newName
        }
        $shp.ConvertToInlineShape() | Out-Null
    }
}

for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $sec = $d.Sections.Item($si)

    for ($hi = 1; $hi -le 3; $hi++) {
        $h = $sec.Headers.Item($hi)
        if ($h.Exists) {
            Rename-InlineShapes $h "image1.jpg" "image2.jpg"
        }
    }

    for ($fi = 1; $fi -le 3; $fi++) {
        $f = $sec.Footers.Item($fi)
        if ($f.Exists) {
            Rename-InlineShapes $f "image2.png" "image1.png"
        }
    }
}
